$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "txtUsername"
$ws.Range("C4").Value = "btnLogin"
$ws.Range("B4").Value = "by_"
$ws.Range("C5").Value = "btnLogin"

$ws.Range("C7").Select()
